$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second worker's row (MADELEY DE JESUS LEON ANTEQUERA, doc 1047482793).
# Deleting the row shifts rows 18-23 up by one and Excel also prunes the now-unused
# shared strings ("1047482793", "MADELEY DE JESUS LEON ANTEQUERA"), which re-aligns
# every other <v> string index automatically.
$ws.Rows("17").Delete()

# Statement date (row 11) moves from 64568 to 45368 - matches the worker row's date.
$ws.Range("E11").Value = 45368

# Only one worker remains now.
$ws.Range("C13").Value = 1
